# Generate Report for Handoff
# Adds a new tracked file (a274dffa-5fa6-4938-8497-f60f5eaeb8b5.md) as row 9
# on the "Overview", "zh-cn" and "de-de" sheets, mirroring the layout of the
# existing rows (1-8) on each sheet.

$wb = $excel.ActiveWorkbook

$hyperFontColor = 15570276   # matches the workbook's custom HyperLink font (RGB 6495ED)
$dateFormat = "yyyy-mm-dd HH:mm:ss"

$uuid = "a274dffa-5fa6-4938-8497-f60f5eaeb8b5"
$mdName = "$uuid.md"
$xlfHash = "933a6e53572780d01ec55f54f68378d2d00278ec"
$xlfZhCn = "$uuid.$xlfHash.zh-cn.xlf"
$xlfDeDe = "$uuid.$xlfHash.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/9b1f53948dc8be0880b89402af7df3e7d74ad5ca/e2e/$mdName"
$xlfZhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea447054e84a72c118c06c98cbb99436612ba7f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhCn"
$xlfDeDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd74924bc305e67677a43b5624affae599a3975c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeDe"

# ---------------------------------------------------------------------------
# Overview sheet: columns A (File Name), B (zh-cn status), C (de-de status),
# D (Latest Handoff Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A9"), $mdUrl, "", "", $mdName) | Out-Null
$wsOverview.Range("A9").Font.Underline = 1
$wsOverview.Range("A9").Font.Color = $hyperFontColor

$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"

$wsOverview.Range("D9").Value = "2016-03-22 06:44:40"
$wsOverview.Range("D9").NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $mdUrl, "", "", $mdName) | Out-Null
$wsZhCn.Range("A9").Font.Underline = 1
$wsZhCn.Range("A9").Font.Color = $hyperFontColor

$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D9"), $xlfZhCnUrl, "", "", $xlfZhCn) | Out-Null
$wsZhCn.Range("D9").Font.Underline = 1
$wsZhCn.Range("D9").Font.Color = $hyperFontColor

$wsZhCn.Range("E9").Value = "2016-03-22 06:44:36"
$wsZhCn.Range("E9").NumberFormat = $dateFormat

$wsZhCn.Range("H9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H9").NumberFormat = $dateFormat

$wsZhCn.Range("J9").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $mdUrl, "", "", $mdName) | Out-Null
$wsDeDe.Range("A9").Font.Underline = 1
$wsDeDe.Range("A9").Font.Color = $hyperFontColor

$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D9"), $xlfDeDeUrl, "", "", $xlfDeDe) | Out-Null
$wsDeDe.Range("D9").Font.Underline = 1
$wsDeDe.Range("D9").Font.Color = $hyperFontColor

$wsDeDe.Range("E9").Value = "2016-03-22 06:44:40"
$wsDeDe.Range("E9").NumberFormat = $dateFormat

$wsDeDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H9").NumberFormat = $dateFormat

$wsDeDe.Range("J9").Value = "Include"

Write-Host "Added handoff row for $mdName to Overview, zh-cn and de-de sheets."
